# Update the cached Power Query results on the Trials Setup sheet:
#   - ROSETTA-Breast-01 (BNT327-05): Progress 0 -> 12.5
#   - REJOICE (MK-5909-003): Days remaining 1 -> 0
#   - REMASTER (CLOU): Days remaining 21 -> 20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 12.5
$ws.Range("B6").Value = 0
$ws.Range("B8").Value = 20
